$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Forces a numeric-looking string to be written as Excel *text* (t="str")
    # instead of being auto-coerced to a number, while leaving the cell's
    # final number format/style back at the workbook default ("General",
    # style index 0) so no stray style survives into the saved file.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

function Set-EmptyTextValue($addr) {
    # Writes an empty *text* cell (t="str" with an empty <v/>), matching
    # what Excel stores when a text cell is cleared via the quote-prefix
    # ('). A plain Value = "" instead deletes the cell outright.
    $r = $ws.Range($addr)
    $r.Value = "'"
    $r.Style = "Normal"
}

# Row 8
$ws.Range("C8").Value = 93

# Row 9
$ws.Range("C9").Value = 80

# Row 10
Set-EmptyTextValue "A10"
$ws.Range("C10").Value = 58
Set-TextValue "D10" "8"
$ws.Range("E10").Value = "Total"
$ws.Range("F10").Value = 0
Set-TextValue "G10" "0.00"

# Row 11
$ws.Range("A11").Value = "%"
$ws.Range("C11").Value = 66
Set-TextValue "D11" "9"
$ws.Range("E11").Value = "Add Tender Premium "
$ws.Range("F11").Value = 0
Set-TextValue "G11" "0.00"

# Row 12
Set-EmptyTextValue "A12"
$ws.Range("C12").Value = 47
Set-TextValue "D12" "10"
$ws.Range("E12").Value = "Grand Total"

# Row 14
Set-TextValue "G14" "0.00"
Set-TextValue "H14" "0.00"

# Row 16
Set-TextValue "G16" "0.00"
Set-TextValue "H16" "0.00"
